$d = $word.ActiveDocument

# --- Update the date line (outside the table) ---
$d.Content.Find.Execute("2025-04-01 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-04-02 Wednesday", 2)

# --- Update the answer table, cell by cell (row, column) ---
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "60÷9=6, 6"
$t.Cell(1,2).Range.Text = "22÷9=2, 4"
$t.Cell(1,3).Range.Text = "46÷4=11, 2"
$t.Cell(1,4).Range.Text = "31÷6=5, 1"
$t.Cell(1,5).Range.Text = "66÷9=7, 3"

$t.Cell(5,1).Range.Text = "39÷9=4, 3"
$t.Cell(5,2).Range.Text = "69÷4=17, 1"
$t.Cell(5,3).Range.Text = "71÷3=23, 2"
$t.Cell(5,4).Range.Text = "61÷9=6, 7"
$t.Cell(5,5).Range.Text = "93÷8=11, 5"

$t.Cell(9,1).Range.Text = "59÷3=19, 2"
$t.Cell(9,2).Range.Text = "90÷2=45, 0"
$t.Cell(9,3).Range.Text = "54÷4=13, 2"
$t.Cell(9,4).Range.Text = "17÷7=2, 3"
$t.Cell(9,5).Range.Text = "11÷2=5, 1"

$t.Cell(13,1).Range.Text = "19÷5=3, 4"
$t.Cell(13,2).Range.Text = "76÷9=8, 4"
$t.Cell(13,3).Range.Text = "57÷3=19, 0"
$t.Cell(13,4).Range.Text = "52÷4=13, 0"
$t.Cell(13,5).Range.Text = "45÷8=5, 5"

$t.Cell(17,1).Range.Text = "84÷2=42, 0"
$t.Cell(17,2).Range.Text = "91÷7=13, 0"
$t.Cell(17,3).Range.Text = "86÷5=17, 1"
$t.Cell(17,4).Range.Text = "76÷4=19, 0"
$t.Cell(17,5).Range.Text = "58÷5=11, 3"

Write-Output "Done"
